$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.630.14'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '1.963.27'
$ws.Range("E3").Value = '  +1.07%  '
$ws.Range("E4").Value = '  -0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '244.44'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.64%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '59.12'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.71%  '
$ws.Range("E8").Value = '  -0.02%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.375'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +3.09%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0812'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -2.11%  '
$ws.Range("E11").Value = '  +0.33%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '22.32'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +3.57%  '
$ws.Range("D13").Value = '2.252.61'
$ws.Range("E13").Value = '  +0.95%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.827'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.08%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '13.74'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("E16").Value = '  +0.75%  '
$ws.Range("D17").Value = '1.966.91'
$ws.Range("E17").Value = '  +1.66%  '
$ws.Range("D18").Value = '36.494.32'
$ws.Range("E18").Value = '  +0.34%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '69.93'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("D20").Value = '0.0₃0858'
$ws.Range("E20").Value = '  -0.08%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '229.01'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("E22").Value = '  +0.74%  '
$ws.Range("E23").Value = '  -0.12%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.44'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.62%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.35'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +3.15%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.141'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +8.71%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.21'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.36%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '160.25'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.91%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '19.46'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.25%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.119'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.86%  '
$ws.Range("E31").Value = '  +0.95%  '
$ws.Range("E32").Value = '  +1.56%  '
$ws.Range("E33").Value = '  -0.99%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.28'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +0.78%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("E36").Value = '  +6.35%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '5.96'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -4.67%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.34'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +11.35%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("E40").Value = '  +1.07%  '
$ws.Range("E41").Value = '  +1.34%  '
$ws.Range("E42").Value = '  +0.24%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.0211'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.27%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '16.07'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.74%  '
$ws.Range("D45").Value = '1.362.45'
$ws.Range("E45").Value = '  +1.15%  '
$ws.Range("E46").Value = '  +1.01%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '87.80'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.37%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.14'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.63%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.83'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.97%  '
$ws.Range("D50").Value = '2.142.20'
$ws.Range("E50").Value = '  +1.00%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '43.83'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -3.06%  '
